# Commit: "manaual discount test data update and keywords"
#
# Updates the Promo test-data workbook:
#   - serial_key column (B2:B63): 307220324rGH -> 307220324WWP
#   - username_admin column (E2:E63): AnandArya -> i9Qa_user1
#   - moves the sheet selection/cursor to L5 (scrolled view towards column F)
#
# (The absPath recorded by Excel under mc:AlternateContent and the
# incidental de-duplication of a couple of unused/duplicate style and
# border table entries are byproducts of the author's local Excel/save
# environment - they have no settable equivalent in the Excel object
# model and no visible effect on the sheet, so they are not reproduced
# here.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$firstRow = 2
$lastRow  = 63

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 2).Value = "307220324WWP"   # B column - serial_key
    $ws.Cells.Item($row, 5).Value = "i9Qa_user1"      # E column - username_admin
}

# Move the active selection (matches the committed sheetView selection change)
$ws.Range("L5").Select()
